$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at the top of the data block (before old row 65),
# pushing the existing rows (old 65..152) down to 68..155.
$ws.Rows("65:67").Insert()

# Seed the 3 new rows with a copy of the (now shifted) rows 68:70 so that the
# constant/common columns (A-C, E-J, Q) and formatting match the rest of the
# table; the varying columns are overwritten individually below.
$ws.Range("A65:T67").Value2 = $ws.Range("A68:T70").Value2

# Row 65 - new weekly entry
$ws.Range("D65").Value2 = 45118
$ws.Range("K65").Value2 = "Fukumoto"
$ws.Range("L65").Value2 = "Segunda"
$ws.Range("M65").Value2 = 300
$ws.Range("N65").Value2 = 650
$ws.Range("O65").Value2 = 700
$ws.Range("P65").Value2 = 675
$ws.Range("R65").Value2 = "Provincia de Melipilla"
$ws.Range("S65").Value2 = 675
$ws.Range("T65").Value2 = 1

# Row 66 - new weekly entry
$ws.Range("D66").Value2 = 45118
$ws.Range("K66").Value2 = "New Hall"
$ws.Range("L66").Value2 = "Segunda"
$ws.Range("M66").Value2 = 270
$ws.Range("N66").Value2 = 750
$ws.Range("O66").Value2 = 800
$ws.Range("P66").Value2 = 775
$ws.Range("R66").Value2 = "Región de O'Higgins"
$ws.Range("S66").Value2 = 775
$ws.Range("T66").Value2 = 1

# Row 67 - new weekly entry
$ws.Range("D67").Value2 = 45118
$ws.Range("K67").Value2 = "New Hall"
$ws.Range("L67").Value2 = "Tercera"
$ws.Range("M67").Value2 = 300
$ws.Range("N67").Value2 = 600
$ws.Range("O67").Value2 = 650
$ws.Range("P67").Value2 = 625
$ws.Range("R67").Value2 = "Región de O'Higgins"
$ws.Range("S67").Value2 = 625
$ws.Range("T67").Value2 = 1
